$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getRelationById")
$ws.Activate()

# Preserve the text of the existing threaded comment currently anchored at C1
# (it will end up on D1 once the new column is inserted before it).
$existingCommentText = $ws.Range("C1").Comment.Text()
$ws.Range("C1").Comment.Delete()

# Insert a new column before the existing "relationId" column (C), shifting
# relationId (and everything to its right) one column to the right.
$ws.Columns("C").Insert()

# New column C picked up the default width; give it the same width used by
# columns A:B (the "wide" style used for the other label/description columns).
$ws.Columns("C").ColumnWidth = 29.67

# Populate the new "label" column.
$ws.Range("C1").Value2 = "label"
$ws.Range("C2").Value2 = "unit"
$ws.Range("C3").Value2 = "is_instance_of"
# C4/C5 intentionally left blank (matches the two "bad request" rows).

# Re-create the comments as threaded comments: the new "label" header comment
# on C1, and the original relationId comment moved along with its cell to D1.
$ws.Range("C1").AddCommentThreaded("If this field is set, test case will use its value to get a list of relations, then pick up the id of the 1st relation as the relationId value to be test and the input parameter ""relationId"" will be ignored.")
$ws.Range("D1").AddCommentThreaded($existingCommentText)

# Match the saved selection/active cell.
$ws.Range("E8").Select()
